$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the trailing column (old BA) and trailing rows (23:24)
$ws.Columns("BA").Delete()
$ws.Rows("23:24").Delete()

# Step 2: update the shifted header/data values row by row
$row1 = @(39583, 39765, 39948, 40130, 40310, 40494, 40676, 40862, 41044, 41228, 41409, 41592, 41774, 41957, 42137, 42321, 42503, 42689, 42867, 43053, 43145, 43235, 43326, 43418, 43510, 43600, 43691, 43783, 43875, 43966, 44068, 44159, 44251, 44341, 44432, 44525, 44617, 44706, 44798, 44890, 44981, 45071, 45163, 45254, 45345, 45436, 45534, 45618, 45713, 45800, 45891)
for ($i = 0; $i -lt $row1.Length; $i++) {
    $v = $row1[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(1, 2+$i).ClearContents()
    } else {
        $ws.Cells.Item(1, 2+$i).Value = $v
    }
}

$row3 = @(5.515135932605308, 2.741420124134053, -2.764179109379705)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $v = $row3[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(3, 2+$i).ClearContents()
    } else {
        $ws.Cells.Item(3, 2+$i).Value = $v
    }
}

$row4 = @(5.566598885960738, 3.4712032095223, -0.2596132895792413, -1.796159642965267, 3.38265053313096)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $v = $row4[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(4, 2+$i).ClearContents()
    } else {
        $ws.Cells.Item(4, 2+$i).Value = $v
    }
}

$row5 = @($null, -0.06191106053077267, -1.072567702512639, 1.421687849828523, 3.372216600811506, 5.813045170083808)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $v = $row5[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(5, 3+$i).ClearContents()
    } else {
        $ws.Cells.Item(5, 3+$i).Value = $v
    }
}

$row6 = @($null, 1.047976847023846, 1.316610894175119, 2.371160938652705, 5.794553957309168, 4.597183386292891)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $v = $row6[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(6, 5+$i).ClearContents()
    } else {
        $ws.Cells.Item(6, 5+$i).Value = $v
    }
}

$row7 = @($null, 2.253457754589228, 5.18838101923762, 6.213583554874536, 4.453225076234824, 2.113053977048707)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $v = $row7[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(7, 7+$i).ClearContents()
    } else {
        $ws.Cells.Item(7, 7+$i).Value = $v
    }
}

$row8 = @($null, 6.4658713724427, 5.226878449528161, 3.80132608443593, 4.443665867214519, 5.262295419893648)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $v = $row8[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(8, 9+$i).ClearContents()
    } else {
        $ws.Cells.Item(8, 9+$i).Value = $v
    }
}

$row9 = @($null, 3.947168334316919, 3.983760133672121, 3.911207409579687, 4.93942156312468, 4.78031692483154)
for ($i = 0; $i -lt $row9.Length; $i++) {
    $v = $row9[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(9, 11+$i).ClearContents()
    } else {
        $ws.Cells.Item(9, 11+$i).Value = $v
    }
}

$row10 = @($null, 3.718405228598609, 4.397175883931581, 4.481753591536197, 5.304093839766777, 5.386945580119185)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $v = $row10[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(10, 13+$i).ClearContents()
    } else {
        $ws.Cells.Item(10, 13+$i).Value = $v
    }
}

$row11 = @($null, 4.481643758166953, 5.391718764442932, 5.291601650505706, 4.773858406282372, 4.884223728030879)
for ($i = 0; $i -lt $row11.Length; $i++) {
    $v = $row11[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(11, 15+$i).ClearContents()
    } else {
        $ws.Cells.Item(11, 15+$i).Value = $v
    }
}

$row12 = @($null, 5.256062028567232, 4.867860541745572, 4.848884192354119, 5.448823958443616, 6.111775080355519, 6.072005530313129, 5.969015988310433)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $v = $row12[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(12, 17+$i).ClearContents()
    } else {
        $ws.Cells.Item(12, 17+$i).Value = $v
    }
}

$row13 = @($null, $null, 4.819877891111335, 5.22223210260726, 5.5540247491221, 5.47980442665561, 5.321568540215016, 5.03659417252571, 3.94180168325462, 4.927928448556984, 3.961704997901161)
for ($i = 0; $i -lt $row13.Length; $i++) {
    $v = $row13[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(13, 18+$i).ClearContents()
    } else {
        $ws.Cells.Item(13, 18+$i).Value = $v
    }
}

$row14 = @($null, $null, $null, 5.484025951319293, 5.40641359612255, 5.259836672320484, 5.077291536394846, 5.802765867180804, 4.342221082693465, 4.009670676786059, 3.818894565497888, 2.502519143054571, -3.840397826549158)
for ($i = 0; $i -lt $row14.Length; $i++) {
    $v = $row14[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(14, 20+$i).ClearContents()
    } else {
        $ws.Cells.Item(14, 20+$i).Value = $v
    }
}

$row15 = @($null, $null, $null, $null, $null, 5.537148498057909, 4.858483928622914, 4.692141525353577, 4.451460339319802, 3.379658261193086, -2.491095770678031, -0.236010050592228, -1.162988086281536, 0.3986977119751156, 0.4839811651348835)
for ($i = 0; $i -lt $row15.Length; $i++) {
    $v = $row15[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(15, 22+$i).ClearContents()
    } else {
        $ws.Cells.Item(15, 22+$i).Value = $v
    }
}

$row16 = @($null, $null, $null, $null, $null, $null, 3.725004516987496, -0.1871946888853837, 0.1963989771473607, -0.8351862827690737, 1.500453746466346, 1.25489130894445, 2.177145583294293, 0.9587999512773893, 1.97557360987699, 2.06342951900429)
for ($i = 0; $i -lt $row16.Length; $i++) {
    $v = $row16[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(16, 25+$i).ClearContents()
    } else {
        $ws.Cells.Item(16, 25+$i).Value = $v
    }
}

$row17 = @($null, $null, $null, $null, $null, -0.6271380860239262, 0.4027662762101869, -0.3850568868524862, 0.0360848986632023, -1.150174078429844, 0.05919149746531627, 0.47378583592943, 0.3896432785800652, -2.242084520390608, -2.311523918755531, -2.156362896191677)
for ($i = 0; $i -lt $row17.Length; $i++) {
    $v = $row17[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(17, 29+$i).ClearContents()
    } else {
        $ws.Cells.Item(17, 29+$i).Value = $v
    }
}

$row18 = @($null, $null, $null, $null, $null, -0.8630989168629921, -0.5798629304026526, -0.3296074516580316, 4.745563085952331, 0.4860776313184267, -0.1671238976421296, 0.7251362782769055, -1.437335768580206, -1.379035600217182, -0.8195740704359578, -0.8205034771073372)
for ($i = 0; $i -lt $row18.Length; $i++) {
    $v = $row18[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(18, 33+$i).ClearContents()
    } else {
        $ws.Cells.Item(18, 33+$i).Value = $v
    }
}

$row19 = @($null, $null, $null, $null, $null, 0.9066518070286467, 0.1286052726408693, 0.7930386201910888, -0.8900729308547461, -1.053916016632561, -0.4351858173977874, -0.1639718104104904, 0.07297157746815053, 0.2518138058242414, 0.4857399523052974, 0.5750555200350504)
for ($i = 0; $i -lt $row19.Length; $i++) {
    $v = $row19[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(19, 37+$i).ClearContents()
    } else {
        $ws.Cells.Item(19, 37+$i).Value = $v
    }
}

$row20 = @($null, $null, $null, $null, $null, -1.099759890829544, -1.022802765586617, -0.5861676960653317, -1.165535347616076, -1.185914719184011, -0.3387208997876479, 0.3466793681147884)
for ($i = 0; $i -lt $row20.Length; $i++) {
    $v = $row20[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(20, 41+$i).ClearContents()
    } else {
        $ws.Cells.Item(20, 41+$i).Value = $v
    }
}

$row21 = @($null, $null, $null, $null, $null, -1.256374469693167, -0.4360715507411128, 0.2338620606656017)
for ($i = 0; $i -lt $row21.Length; $i++) {
    $v = $row21[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(21, 45+$i).ClearContents()
    } else {
        $ws.Cells.Item(21, 45+$i).Value = $v
    }
}

$row22 = @($null, $null, $null, $null)
for ($i = 0; $i -lt $row22.Length; $i++) {
    $v = $row22[$i]
    if ($null -eq $v) {
        $ws.Cells.Item(22, 49+$i).ClearContents()
    } else {
        $ws.Cells.Item(22, 49+$i).Value = $v
    }
}
